$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 699, shifting existing rows 699-740 down to 700-741.
$ws.Rows(699).Insert()

# Populate the newly inserted row 699 with the new data point.
$ws.Range("B699").Value = "月"
$ws.Range("C699").Value = 5
$ws.Range("D699").Value = 19

# Column A holds a date-formatted string ("2026/01/26"). Excel's COM layer
# auto-parses bare yyyy/mm/dd text into a date serial, so force literal text
# via the quote-prefix convention, then strip the resulting style back to
# the default so no stray formatting is introduced.
$ws.Range("A699").Value = "'2026/01/26"
$ws.Range("A699").Style = "Normal"
